# Update database and change read_price algorithm
# Shifts the 5 displayed financial periods one column to the left (dropping the
# oldest period) and appends the newest period (1401/12) in column H, along
# with its related publish-date label. Also converts the D15 cell that used to
# hold a literal "-" placeholder into a real numeric value (the old E15 value),
# matching the shift-left pattern applied to the rest of the data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header row 8: financial period labels (columns D:H)
# ---------------------------------------------------------------------------
$ws.Range("D8").Value = "12 ماهه منتهی به 1397/12"
$ws.Range("E8").Value = "12 ماهه منتهی به 1398/12"
$ws.Range("F8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("G8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("H8").Value = "12 ماهه منتهی به 1401/12"

# ---------------------------------------------------------------------------
# Header row 9: publish-date labels (columns D:H)
# ---------------------------------------------------------------------------
$ws.Range("D9").Value = "1399-04-08 (11)"
$ws.Range("E9").Value = "1400-04-06 (11)"
$ws.Range("F9").Value = "1401-02-21 (10)"
$ws.Range("G9").Value = "1402-02-23 (10)"
$ws.Range("H9").Value = "1402-02-23 (2)"

# ---------------------------------------------------------------------------
# Row 11: فروش (Sales)
# ---------------------------------------------------------------------------
$ws.Range("D11").Value = 4611283
$ws.Range("E11").Value = 5141567
$ws.Range("F11").Value = 9013378
$ws.Range("G11").Value = 27693494
$ws.Range("H11").Value = 91661180

# ---------------------------------------------------------------------------
# Row 12: بهای تمام شده کالای فروش رفته (Cost of goods sold)
# ---------------------------------------------------------------------------
$ws.Range("D12").Value = -3037243
$ws.Range("E12").Value = -1882242
$ws.Range("F12").Value = -3334760
$ws.Range("G12").Value = -13034835
$ws.Range("H12").Value = -55297002

# ---------------------------------------------------------------------------
# Row 13: سود (زیان) ناخالص (Gross profit)
# ---------------------------------------------------------------------------
$ws.Range("D13").Value = 1574040
$ws.Range("E13").Value = 3259325
$ws.Range("F13").Value = 5678618
$ws.Range("G13").Value = 14658659
$ws.Range("H13").Value = 36364178

# ---------------------------------------------------------------------------
# Row 14: هزینه های عمومی, اداری و تشکیلاتی (G&A expenses)
# ---------------------------------------------------------------------------
$ws.Range("D14").Value = -494539
$ws.Range("E14").Value = -522730
$ws.Range("F14").Value = -897281
$ws.Range("G14").Value = -2088362
$ws.Range("H14").Value = -9847885

# ---------------------------------------------------------------------------
# Row 15: هزینه کاهش ارزش دریافتنی‌ها (هزینه استثنایی)
# D15 used to hold the literal text "-"; it now carries the numeric value
# that used to sit in E15, continuing the shift-left pattern.
# ---------------------------------------------------------------------------
$ws.Range("D15").Value = -79368
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 0

# ---------------------------------------------------------------------------
# Row 16: خالص سایر درامدها (هزینه ها) ی عملیاتی
# ---------------------------------------------------------------------------
$ws.Range("D16").Value = 53604
$ws.Range("E16").Value = -39716
$ws.Range("F16").Value = 98273
$ws.Range("G16").Value = 217344
$ws.Range("H16").Value = 201919

# ---------------------------------------------------------------------------
# Row 17: سود (زیان) عملیاتی (Operating profit/loss)
# ---------------------------------------------------------------------------
$ws.Range("D17").Value = 1053737
$ws.Range("E17").Value = 2696879
$ws.Range("F17").Value = 4879610
$ws.Range("G17").Value = 12787641
$ws.Range("H17").Value = 26718212

# ---------------------------------------------------------------------------
# Row 18: هزینه های مالی (Financial expenses)
# ---------------------------------------------------------------------------
$ws.Range("D18").Value = -328884
$ws.Range("E18").Value = -129619
$ws.Range("F18").Value = -321760
$ws.Range("G18").Value = -1953699
$ws.Range("H18").Value = -5242657

# ---------------------------------------------------------------------------
# Row 19: خالص سایر درامدها و هزینه های غیرعملیاتی
# ---------------------------------------------------------------------------
$ws.Range("D19").Value = 13173
$ws.Range("E19").Value = 133010
$ws.Range("F19").Value = 177309
$ws.Range("G19").Value = 715341
$ws.Range("H19").Value = 403406

# ---------------------------------------------------------------------------
# Row 20: سود (زیان) خالص عملیات در حال تداوم قبل از مالیات
# ---------------------------------------------------------------------------
$ws.Range("D20").Value = 738026
$ws.Range("E20").Value = 2700270
$ws.Range("F20").Value = 4735159
$ws.Range("G20").Value = 11549283
$ws.Range("H20").Value = 21878961

# ---------------------------------------------------------------------------
# Row 21: مالیات (Tax)
# ---------------------------------------------------------------------------
$ws.Range("D21").Value = -172231
$ws.Range("E21").Value = -506585
$ws.Range("F21").Value = -598746
$ws.Range("G21").Value = -1562392
$ws.Range("H21").Value = -4312147

# ---------------------------------------------------------------------------
# Row 22: سود (زیان) خالص عملیات در حال تداوم
# ---------------------------------------------------------------------------
$ws.Range("D22").Value = 565795
$ws.Range("E22").Value = 2193685
$ws.Range("F22").Value = 4136413
$ws.Range("G22").Value = 9986891
$ws.Range("H22").Value = 17566814

# ---------------------------------------------------------------------------
# Row 23: سود (زیان) عملیات متوقف شده پس از اثر مالیاتی (unchanged, all zero)
# ---------------------------------------------------------------------------
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 0

# ---------------------------------------------------------------------------
# Row 24: سود (زیان) خالص (Net profit/loss)
# ---------------------------------------------------------------------------
$ws.Range("D24").Value = 565795
$ws.Range("E24").Value = 2193685
$ws.Range("F24").Value = 4136413
$ws.Range("G24").Value = 9986891
$ws.Range("H24").Value = 17566814

# ---------------------------------------------------------------------------
# Row 25: سود هر سهم پس از کسر مالیات (EPS after tax)
# ---------------------------------------------------------------------------
$ws.Range("D25").Value = 566
$ws.Range("E25").Value = 2194
$ws.Range("F25").Value = 4136
$ws.Range("G25").Value = 591
$ws.Range("H25").Value = 1039

# ---------------------------------------------------------------------------
# Row 26: سرمایه (Capital)
# ---------------------------------------------------------------------------
$ws.Range("D26").Value = 1000000
$ws.Range("E26").Value = 1000000
$ws.Range("F26").Value = 1000000
$ws.Range("G26").Value = 16900000
$ws.Range("H26").Value = 16900000

# ---------------------------------------------------------------------------
# Row 27: سود هر سهم بر اساس آخرین سرمایه
# ---------------------------------------------------------------------------
$ws.Range("D27").Value = 33
$ws.Range("E27").Value = 130
$ws.Range("F27").Value = 245
$ws.Range("G27").Value = 591
$ws.Range("H27").Value = 1039
